# Aula_02_BD.pptx — "Modelos conceitual, lógico e físico finalizados"
#
# 1) The auto-updating "datetimeFigureOut" date placeholders (slide master,
#    all 12 slide layouts and the handout master) get their cached date
#    text refreshed from 30/07/2021 to 03/08/2021.
# 2) Slide 12's "SQL Server" history paragraph is corrected:
#    "em 1984." -> "em 1988."

$p = $ppt.ActivePresentation

$oldDate = "30/07/2021"
$newDate = "03/08/2021"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide master's own date placeholder.
Update-DatePlaceholder($p.SlideMaster.Shapes)

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder($layouts.Item($j).Shapes)
}

# Handout master's date placeholder.
Update-DatePlaceholder($p.HandoutMaster.Shapes)

# Slide 12: "Foi criado em parceria entre Microsoft e a Sybase, em 1984." -> "... em 1988."
$slide = $p.Slides.Item(12)
$shape = $slide.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$fullText = $tr.Text
$target = "em 1984."
$idx = $fullText.IndexOf($target)
if ($idx -ge 0) {
    $chars = $tr.Characters($idx + 1, $target.Length)
    $chars.Text = "em 1988."
}
